# Fruta / hortaliza, semanal
# Insert a new week's worth of price rows (date 44504) into the Cebollín
# subset sheet, just before the existing row 566 block, pushing the rest
# of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 566-567; everything from the old row 566
# onward shifts down to 568 onward.
$ws.Rows("566:567").Insert()

# New row 566: Primera grade, week of 44504 (2021-11-26)
$ws.Range("A566").Value = 9
$ws.Range("B566").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C566").Value = "Metropolitana"
$ws.Range("D566").Value = 44504
$ws.Range("E566").Value = 13
$ws.Range("F566").Value = 100112037
$ws.Range("G566").Value = "Cebollín"
$ws.Range("H566").Value = "Sin especificar"
$ws.Range("I566").Value = "Primera"
$ws.Range("J566").Value = 250
$ws.Range("K566").Value = 2200
$ws.Range("L566").Value = 2500
$ws.Range("M566").Value = 2350
$ws.Range("N566").Value = "$/paquete 36 unidades"
$ws.Range("O566").Value = "Región Metropolitana"
$ws.Range("P566").Value = 65
$ws.Range("Q566").Value = 36
$ws.Range("R566").Value = "Hortaliza"

# New row 567: Segunda grade, week of 44504 (2021-11-26)
$ws.Range("A567").Value = 9
$ws.Range("B567").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C567").Value = "Metropolitana"
$ws.Range("D567").Value = 44504
$ws.Range("E567").Value = 13
$ws.Range("F567").Value = 100112037
$ws.Range("G567").Value = "Cebollín"
$ws.Range("H567").Value = "Sin especificar"
$ws.Range("I567").Value = "Segunda"
$ws.Range("J567").Value = 106
$ws.Range("K567").Value = 1700
$ws.Range("L567").Value = 2000
$ws.Range("M567").Value = 1850
$ws.Range("N567").Value = "$/paquete 36 unidades"
$ws.Range("O567").Value = "Región Metropolitana"
$ws.Range("P567").Value = 51
$ws.Range("Q567").Value = 36
$ws.Range("R567").Value = "Hortaliza"
